$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.926126621258732
$ws.Range("D2").Value = 3.338002904462885
$ws.Range("E2").Value = 16.66695836158539
$ws.Range("F2").Value = 27.18703899061369
$ws.Range("G2").Value = 3.567828460017615
$ws.Range("O2").Value = 23.19951237517529

$ws.Range("C3").Value = 4.750213261576099
$ws.Range("D3").Value = 3.332266480516632
$ws.Range("E3").Value = 15.70294774289622
$ws.Range("F3").Value = 26.31248136070441
$ws.Range("G3").Value = 3.572303459150941
$ws.Range("O3").Value = 22.57619591036849

$ws.Range("C4").Value = 4.640614124971141
$ws.Range("D4").Value = 3.329717909951419
$ws.Range("E4").Value = 15.08611945890479
$ws.Range("F4").Value = 25.77068349195109
$ws.Range("G4").Value = 3.575189172148516
$ws.Range("O4").Value = 22.19320626190192

$ws.Range("C5").Value = 4.595631543631719
$ws.Range("D5").Value = 3.328923112997732
$ws.Range("E5").Value = 14.82875600235974
$ws.Range("F5").Value = 25.54907308206905
$ws.Range("G5").Value = 3.576399996207762
$ws.Range("O5").Value = 22.0373268116054

$ws.Range("C6").Value = 4.588145329460033
$ws.Range("D6").Value = 3.328805825206065
$ws.Range("E6").Value = 14.78566738895243
$ws.Range("F6").Value = 25.51223691030425
$ws.Range("G6").Value = 3.576603163521204
$ws.Range("O6").Value = 22.01146251924012

$ws.Range("C7").Value = 4.64000866065126
$ws.Range("D7").Value = 3.329706205593732
$ws.Range("E7").Value = 15.08267247213794
$ws.Range("F7").Value = 25.76769758027978
$ws.Range("G7").Value = 3.57520536033398
$ws.Range("O7").Value = 22.19110288947422

$ws.Range("C8").Value = 4.86584670503044
$ws.Range("D8").Value = 3.335822339363518
$ws.Range("E8").Value = 16.33988363937718
$ws.Range("F8").Value = 26.88671595412688
$ws.Range("G8").Value = 3.569342884976451
$ws.Range("O8").Value = 22.98480201195886

$ws.Range("C9").Value = 5.292972261085125
$ws.Range("D9").Value = 3.355579124847962
$ws.Range("E9").Value = 18.74361448004077
$ws.Range("F9").Value = 29.02717102200161
$ws.Range("G9").Value = 3.558934611588261
$ws.Range("O9").Value = 24.52821902330734

$ws.Range("C10").Value = 5.593444883520056
$ws.Range("D10").Value = 3.374865229620673
$ws.Range("E10").Value = 20.43372571101
$ws.Range("F10").Value = 30.54778217139857
$ws.Range("G10").Value = 3.551940698096475
$ws.Range("O10").Value = 25.64059188266678

$ws.Range("C11").Value = 5.726582748968817
$ws.Range("D11").Value = 3.384678399948782
$ws.Range("E11").Value = 21.16060891215215
$ws.Range("F11").Value = 31.22479402351989
$ws.Range("G11").Value = 3.548898584063292
$ws.Range("O11").Value = 26.13936855543927

$ws.Range("C12").Value = 5.776441726129736
$ws.Range("D12").Value = 3.38854379853924
$ws.Range("E12").Value = 21.42987029500259
$ws.Range("F12").Value = 31.47879319919044
$ws.Range("G12").Value = 3.547766495598578
$ws.Range("O12").Value = 26.32700660914464

$ws.Range("C13").Value = 5.765729185756652
$ws.Range("D13").Value = 3.387704677737932
$ws.Range("E13").Value = 21.37214596823448
$ws.Range("F13").Value = 31.42419911024617
$ws.Range("G13").Value = 3.548009428843606
$ws.Range("O13").Value = 26.28665339029794

$ws.Range("C14").Value = 5.730696098716321
$ws.Range("D14").Value = 3.384993416870912
$ws.Range("E14").Value = 21.18288109909068
$ws.Range("F14").Value = 31.24573955303894
$ws.Range("G14").Value = 3.548805048668435
$ws.Range("O14").Value = 26.15483152136525

$ws.Range("C15").Value = 5.709163424527785
$ws.Range("D15").Value = 3.383352134246714
$ws.Range("E15").Value = 21.06617155094796
$ws.Range("F15").Value = 31.13611214935319
$ws.Range("G15").Value = 3.549294974859055
$ws.Range("O15").Value = 26.07392009104961

$ws.Range("C16").Value = 5.584668440015265
$ws.Range("D16").Value = 3.374244831144191
$ws.Range("E16").Value = 20.38538131221744
$ws.Range("F16").Value = 30.50321931614273
$ws.Range("G16").Value = 3.552142299957014
$ws.Range("O16").Value = 25.60783271292768

$ws.Range("C17").Value = 5.507350503590494
$ws.Range("D17").Value = 3.368924076746697
$ws.Range("E17").Value = 19.95702076391891
$ws.Range("F17").Value = 30.11100033487733
$ws.Range("G17").Value = 3.553924645731944
$ws.Range("O17").Value = 25.31990249809807

$ws.Range("C18").Value = 5.462547692818938
$ws.Range("D18").Value = 3.36596161706637
$ws.Range("E18").Value = 19.70668754964538
$ws.Range("F18").Value = 29.88403404494042
$ws.Range("G18").Value = 3.554962939315402
$ws.Range("O18").Value = 25.1536211536642

$ws.Range("C19").Value = 5.44732281414514
$ws.Range("D19").Value = 3.36497539429062
$ws.Range("E19").Value = 19.6212491263944
$ws.Range("F19").Value = 29.80695963907823
$ws.Range("G19").Value = 3.555316748541758
$ws.Range("O19").Value = 25.0972120508129

$ws.Range("C20").Value = 5.515615806980166
$ws.Range("D20").Value = 3.369480348570722
$ws.Range("E20").Value = 20.00302915766162
$ws.Range("F20").Value = 30.15289666682768
$ws.Range("G20").Value = 3.553733553664565
$ws.Range("O20").Value = 25.35062412894841

$ws.Range("C21").Value = 5.741001641805984
$ws.Range("D21").Value = 3.38578572861479
$ws.Range("E21").Value = 21.23863512412189
$ws.Range("F21").Value = 31.29822367207156
$ws.Range("G21").Value = 3.548570817187844
$ws.Range("O21").Value = 26.19358587191

$ws.Range("C22").Value = 5.885037196549966
$ws.Range("D22").Value = 3.397312611908109
$ws.Range("E22").Value = 22.0112587194385
$ws.Range("F22").Value = 32.03284516023584
$ws.Range("G22").Value = 3.54531255446324
$ws.Range("O22").Value = 26.73721587234338

$ws.Range("C23").Value = 5.808475767081984
$ws.Range("D23").Value = 3.39108096922618
$ws.Range("E23").Value = 21.60207667208536
$ws.Range("F23").Value = 31.6421139587905
$ws.Range("G23").Value = 3.547040999811411
$ws.Range("O23").Value = 26.44779810753269

$ws.Range("C24").Value = 5.511880155142912
$ws.Range("D24").Value = 3.369228557453475
$ws.Range("E24").Value = 19.98224142129509
$ws.Range("F24").Value = 30.13395990796414
$ws.Range("G24").Value = 3.55381990401785
$ws.Range("O24").Value = 25.33673718956986

$ws.Range("C25").Value = 5.179520665868943
$ws.Range("D25").Value = 3.34939844804807
$ws.Range("E25").Value = 18.08374502840276
$ws.Range("F25").Value = 28.45597794450965
$ws.Range("G25").Value = 3.561634903240414
$ws.Range("O25").Value = 24.11354218685661
